$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new URL as a hyperlink in cell A3 (this also writes the cell's
# display text / shared string and applies the built-in Hyperlink style).
$ws.Hyperlinks.Add($ws.Range("A3"), "https://www.tech.gov.sg/") | Out-Null

# Match the author's final selection (cell A3 was the active cell on save).
$ws.Range("A3").Select() | Out-Null
